# Weekly update: insert a new price record for "Espinaca" (Agrícola del
# Norte S.A. de Arica) as the newest row, pushing the existing history
# down by one row (row 82 -> 83, ..., row 108 -> 109).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 82; Excel shifts rows 82:108 down to 83:109.
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with this week's record.
$ws.Range("A82").Value = 1
$ws.Range("B82").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C82").Value = "Arica y Parinacota"
$ws.Range("D82").Value = 45120
$ws.Range("E82").Value = 15
$ws.Range("F82").Value = 100112012
$ws.Range("G82").Value = "Espinaca"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 450
$ws.Range("K82").Value = 1400
$ws.Range("L82").Value = 1500
$ws.Range("M82").Value = 1433
$ws.Range("N82").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O82").Value = "Región de Arica y Parinacota"
$ws.Range("P82").Value = 478
$ws.Range("Q82").Value = 3
$ws.Range("R82").Value = "Hortaliza"
